# Add a "Save" column (H) to the s_vals worksheet, mirroring the existing
# header style used by the other header cells (B1:G1) and filling in the
# per-row 0/1 values from column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled like the rest of row 1 (bold, bordered, centered) by
# copying the existing header format from G1 onto the new H1 header.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Data values for the new "Save" column, rows 2-17.
$saveValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}
